$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 23.069913
$ws.Range("H2").Value = 69.20973899999998
$ws.Range("I2").Value = 0.9355059672894461
$ws.Range("J2").Value = 0.9355059672894461
$ws.Range("M2").Value = 0.7285076666666667
$ws.Range("N2").Value = 2.185523
$ws.Range("O2").Value = 0.1384760821597099
$ws.Range("P2").Value = 0.1384760821597099
$ws.Range("Q2").Value = 16.806608489833
$ws.Range("R2").Value = 151.259476408497
$ws.Range("S2").Value = 0.1295452011872723
$ws.Range("T2").Value = 0.1295452011872723
$ws.Range("G3").Value = 23.069913
$ws.Range("H3").Value = 69.20973899999998
$ws.Range("I3").Value = 0.9355059672894461
$ws.Range("J3").Value = 0.9355059672894461
$ws.Range("O3").Value = 0.6813230330092965
$ws.Range("P3").Value = 0.6813230330092966
$ws.Range("Q3").Value = 82.69102716009998
$ws.Range("R3").Value = 744.2192444408998
$ws.Range("S3").Value = 0.6373817630319412
$ws.Range("T3").Value = 0.6373817630319413
$ws.Range("G4").Value = 23.069913
$ws.Range("H4").Value = 69.20973899999998
$ws.Range("I4").Value = 0.9355059672894461
$ws.Range("J4").Value = 0.9355059672894461
$ws.Range("O4").Value = 0.1802008848309935
$ws.Range("P4").Value = 0.1802008848309935
$ws.Range("Q4").Value = 21.870677402492
$ws.Range("R4").Value = 196.836096622428
$ws.Range("S4").Value = 0.1685790030702326
$ws.Range("T4").Value = 0.1685790030702326
$ws.Range("I5").Value = 0.03069514654402774
$ws.Range("J5").Value = 0.03069514654402774
$ws.Range("M5").Value = 0.7285076666666667
$ws.Range("N5").Value = 2.185523
$ws.Range("O5").Value = 0.1384760821597099
$ws.Range("P5").Value = 0.1384760821597099
$ws.Range("Q5").Value = 0.5514463066422222
$ws.Range("R5").Value = 4.963016759779999
$ws.Range("S5").Value = 0.004250543634735122
$ws.Range("T5").Value = 0.004250543634735122
$ws.Range("I6").Value = 0.03069514654402774
$ws.Range("J6").Value = 0.03069514654402774
$ws.Range("O6").Value = 0.6813230330092965
$ws.Range("P6").Value = 0.6813230330092966
$ws.Range("S6").Value = 0.02091331034204181
$ws.Range("T6").Value = 0.02091331034204181
$ws.Range("I7").Value = 0.03069514654402774
$ws.Range("J7").Value = 0.03069514654402774
$ws.Range("O7").Value = 0.1802008848309935
$ws.Range("P7").Value = 0.1802008848309935
$ws.Range("S7").Value = 0.00553129256725081
$ws.Range("T7").Value = 0.005531292567250811
$ws.Range("I8").Value = 0.03379888616652608
$ws.Range("J8").Value = 0.03379888616652608
$ws.Range("M8").Value = 0.7285076666666667
$ws.Range("N8").Value = 2.185523
$ws.Range("O8").Value = 0.1384760821597099
$ws.Range("P8").Value = 0.1384760821597099
$ws.Range("Q8").Value = 0.6072057977771111
$ws.Range("R8").Value = 5.464852179994
$ws.Range("S8").Value = 0.004680337337702549
$ws.Range("T8").Value = 0.004680337337702549
$ws.Range("I9").Value = 0.03379888616652608
$ws.Range("J9").Value = 0.03379888616652608
$ws.Range("O9").Value = 0.6813230330092965
$ws.Range("P9").Value = 0.6813230330092966
$ws.Range("S9").Value = 0.0230279596353135
$ws.Range("T9").Value = 0.02302795963531351
$ws.Range("I10").Value = 0.03379888616652608
$ws.Range("J10").Value = 0.03379888616652608
$ws.Range("O10").Value = 0.1802008848309935
$ws.Range("P10").Value = 0.1802008848309935
$ws.Range("S10").Value = 0.006090589193510024
$ws.Range("T10").Value = 0.006090589193510024
